# Generate Report for Archive
#
# Updates the localization status: the "Ready for handoff" status becomes
# "In Translation" everywhere it is used (Overview!E2/F2, zh-cn!C2,
# de-de!C2), and the now-narrower Status columns are shrunk to fit the
# shorter text (Overview columns E & F, and column C on the zh-cn / de-de
# sheets).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Update the status text wherever it appears.
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Shrink the corresponding columns to match the narrower text.
$newColumnWidth = 12.5

$wsOverview.Range("E1:F1").ColumnWidth = $newColumnWidth
$wsZhCn.Range("C1").ColumnWidth = $newColumnWidth
$wsDeDe.Range("C1").ColumnWidth = $newColumnWidth
